$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 148 (shifts existing rows 148..238 down to 149..239)
$ws.Rows.Item(148).Insert()

# Populate the newly inserted row 148 with its data.
# Columns A, B, C, E, F, G, H, I, N, Q, R keep the same values as the row that was there before
# (Excel's Insert already copies formatting, so we just need to set the values).
$ws.Range("A148").Value = 11
$ws.Range("B148").Value = 'Vega Monumental Concepción'
$ws.Range("C148").Value = 'Bíobío'
$ws.Range("D148").Value = 45097
$ws.Range("E148").Value = 8
$ws.Range("F148").Value = 100112032
$ws.Range("G148").Value = 'Zapallo italiano'
$ws.Range("H148").Value = 'Sin especificar'
$ws.Range("I148").Value = 'Primera'
$ws.Range("J148").Value = 250
$ws.Range("K148").Value = 11000
$ws.Range("L148").Value = 12000
$ws.Range("M148").Value = 11600
$ws.Range("N148").Value = '$/caja 50 unidades'
$ws.Range("O148").Value = 'Región de Arica y Parinacota'
$ws.Range("P148").Value = 232
$ws.Range("Q148").Value = 50
$ws.Range("R148").Value = 'Hortaliza'
